# Backup before dimension reduction
# Shift the "qN" index labels in column A up by one: row r gets value "q" + (r-2)
# i.e. A2: q1 -> q0, A3: q2 -> q1, ..., A97: q96 -> q95

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 97; $r++) {
    $ws.Cells.Item($r, 1).Value = "q" + ($r - 2)
}
